# Update cryptocurrency price and volume(1h) data per GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'308.46"
$ws.Range('E2').Value = "'0.55%"
$ws.Range('D3').Value = "'40.82"
$ws.Range('E3').Value = "'2.41%"
$ws.Range('D4').Value = "'5.117"
$ws.Range('E4').Value = "'-0.21%"
$ws.Range('D5').Value = "'0.07613"
$ws.Range('E5').Value = "'-1.21%"
$ws.Range('D6').Value = "'1.625"
$ws.Range('E6').Value = "'-0.33%"
$ws.Range('D7').Value = "'0.9018"
$ws.Range('E7').Value = "'2.40%"
$ws.Range('E8').Value = "'-0.28%"
$ws.Range('D9').Value = "'0.1098"
$ws.Range('E9').Value = "'9.57%"
$ws.Range('D10').Value = "'0.1773"
$ws.Range('E10').Value = "'1.54%"
$ws.Range('D11').Value = "'0.09190"
$ws.Range('E11').Value = "'3.00%"
$ws.Range('D12').Value = "'0.04176"
$ws.Range('E12').Value = "'-5.25%"
$ws.Range('E13').Value = "'-0.55%"
$ws.Range('D14').Value = "'0.001258"
$ws.Range('E14').Value = "'-0.23%"
$ws.Range('D15').Value = "'0.005862"
$ws.Range('E15').Value = "'-1.04%"
$ws.Range('E16').Value = "'-0.03%"
$ws.Range('D17').Value = "'4.254"
$ws.Range('E17').Value = "'0.37%"
$ws.Range('E18').Value = "'-0.73%"
$ws.Range('D19').Value = "'6.555"
$ws.Range('E19').Value = "'-6.30%"
$ws.Range('D20').Value = "'0.1360"
$ws.Range('E20').Value = "'2.21%"
$ws.Range('E21').Value = "'-10.64%"
$ws.Range('D22').Value = "'0.04066"
$ws.Range('E22').Value = "'-2.18%"
$ws.Range('D23').Value = "'0.001223"
$ws.Range('E23').Value = "'2.41%"
$ws.Range('D24').Value = "'0.004092"
$ws.Range('E24').Value = "'-0.06%"
$ws.Range('D38').Value = "'0.02379"
$ws.Range('E38').Value = "'1.67%"
$ws.Range('D39').Value = "'0.05186"
$ws.Range('E39').Value = "'0.80%"
$ws.Range('D40').Value = "'0.007794"
$ws.Range('E40').Value = "'-1.87%"
$ws.Range('E41').Value = "'-1.68%"
$ws.Range('D42').Value = "'0.006763"
$ws.Range('E42').Value = "'6.74%"
$ws.Range('E43').Value = "'-0.98%"
$ws.Range('D44').Value = "'0.008347"
$ws.Range('E44').Value = "'-2.83%"
$ws.Range('D45').Value = "'0.3077"
$ws.Range('E45').Value = "'0.74%"
$ws.Range('D46').Value = "'0.00006951"
$ws.Range('E46').Value = "'6.91%"
$ws.Range('D47').Value = "'0.00000000751"
$ws.Range('E47').Value = "'0.04%"
$ws.Range('D48').Value = "'0.03181"
$ws.Range('E48').Value = "'398.19%"
$ws.Range('E49').Value = "'-39.98%"
$ws.Range('E50').Value = "'0.04%"
$ws.Range('E51').Value = "'0.04%"
